# Adding two extra columns for capturing the satisfaction level for existing partners
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G and H (Current partner | [NEW] | [NEW] | Prefered Partners ...)
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(8).Insert()

# Headers for the two new columns
$ws.Range("G1").Value = "Current Patner Feedback"
$ws.Range("H1").Value = "Current Patner Feedback Reason"

# New phone number for row 2
$ws.Range("A2").Value = 44444444441

# New data for the new columns
$ws.Range("G2").Value = "Satisfied"
$ws.Range("H2").Value = "Test 1"

$ws.Range("G3").Value = "Dissatisfied"
$ws.Range("H3").Value = "Test 3"

$ws.Range("G4").Value = "Extremely Dissatisfied"

# Column widths for the two new columns
$ws.Columns.Item(7).ColumnWidth = 27.226666666699998
$ws.Columns.Item(8).ColumnWidth = 35.4266666667

# Update selection / view
$ws.Range("N3").Select()
$excel.ActiveWindow.ScrollColumn = 9
